$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing data rows (9 and 10) down into the
# new rows (11 and 12) before filling in values, so number formats /
# borders / fonts / wrap match the rest of the table.
$ws.Range("A9:E9").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A10:E10").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

# Row 11
$ws.Range("A11").Value = 44947
$ws.Range("B11").Value = 0.63055555555555554
$ws.Range("C11").Value = 0.69305555555555554
$ws.Range("D11").Value = "Generated mock data and tested tutor and student availability functionality."
$ws.Range("E11").Value = 1.5

# Row 12
$ws.Range("A12").Value = 44948
$ws.Range("B12").Value = 0.44444444444444442
$ws.Range("C12").Value = 0.4861111111111111
$ws.Range("D12").Value = "Worked on a way of finding tutors available at a specific time and seeing how many times they have worked with a specific student."
$ws.Range("E12").Value = 1

# The descriptions wrap across multiple lines in column D, so the rows
# grow taller to fit the text (matches how Excel auto-sizes these rows).
$ws.Rows.Item(11).RowHeight = 28
$ws.Rows.Item(12).RowHeight = 56

$ws.Range("F13").Select()
